$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workhours")

# The existing "Sum:" separator/formula block currently sits at rows 16-18.
# Insert four fresh rows above it so three new work-log entries (plus one
# blank spacer row, matching the existing table's rhythm) can be recorded
# while pushing the separator/Sum rows down to 20-22.
$ws.Rows("16:19").Insert()

# Copy the formatting (date / duration / description cell styles) from the
# row directly above (row 15) down across the four new rows so they match
# the rest of the log table exactly.
$ws.Range("B15:D15").Copy()
$ws.Range("B16:D19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New work log entries.
$ws.Range("B16").Value = 44947
$ws.Range("C16").Value = 0.041666666666666664
$ws.Range("D16").Value = "Drawing sketches"

$ws.Range("B17").Value = 44948
$ws.Range("C17").Value = 0.041666666666666664
$ws.Range("D17").Value = "Drawing sketches"

$ws.Range("B18").Value = 44949
$ws.Range("C18").Value = 0.16666666666666666
$ws.Range("D18").Value = "Designing electronics and printing first test parts"

# Row 19 stays an empty spacer row (same as row 20 already is), just like
# the original blank row before the Sum line.

[void]$ws.Range("D26").Select()
